$d = $word.ActiveDocument

# --- Table 1 is the "Tidsrapport" (time report) table for this iteration ---
$t = $d.Tables.Item(1)

# Row 4: "Ronnie: Livsystem" -> Status "Pabörjad" -> "Klar"; Verklig tid "5" -> "4"
$t.Cell(4, 3).Range.Text = "Klar"
$t.Cell(4, 5).Range.Text = "4"

# Row 5: "Ronnie: Attackkodning" -> Status "Pabörjad" -> "Klar"
$t.Cell(5, 3).Range.Text = "Klar"

# Row 12: "Övrig dokumentation" -> Verklig tid "1" -> "1,5"
$t.Cell(12, 5).Range.Text = "1,5"

# --- Move the "_GoBack" bookmark from the old edit location (inside the
# "Mål" paragraph) to the end of the text we just typed in row 5 ("Klar"),
# matching Word's behaviour of tracking the most recent edit position. ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$statusCell = $t.Cell(5, 3)
$statusRange = $d.Range($statusCell.Range.Start, $statusCell.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $statusRange)

# --- The footer page-number field's cached display value needs to be
# refreshed now that the document repaginates to a single page. ---
$footer = $d.Sections.Item(1).Footers.Item(1)
$pageField = $footer.Range.Fields.Item(1)
$pageField.Result.Select()
$word.Selection.Text = "1"
